# Updated cryptos list on Sun May 26 04:50:41 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.996.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.33%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.751.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.53%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'602.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.07%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'165.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.61%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.746.07"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.76%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.30%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.173"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.74%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.50%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'37.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.31%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.29%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.381.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.37%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.745.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.50%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'68.998.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.24%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'7.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.25%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'17.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.07%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +5.07%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'490.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.23%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.26%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'84.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.14%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -2.47%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.26%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.79%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.07%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.60%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'8.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.90%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.41%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.900.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'31.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.688.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.35%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.95%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.84%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.02%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +3.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.10%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +7.92%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.82%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'48.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.85%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.19%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'423.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.75%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'8.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.32%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'40.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.52%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'141.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.42%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'ONDO"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'1.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.89%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Maker"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'2.782.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.74%  "
$ws.Range("E51").Style = "Normal"
